$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $row1, $row2, $firstCol, $lastCol) {
    $rng1 = $ws.Range($ws.Cells.Item($row1, $firstCol), $ws.Cells.Item($row1, $lastCol))
    $rng2 = $ws.Range($ws.Cells.Item($row2, $firstCol), $ws.Cells.Item($row2, $lastCol))

    $vals1 = $rng1.Value2
    $vals2 = $rng2.Value2

    $rng1.Value2 = $vals2
    $rng2.Value2 = $vals1
}

# Swap data (columns B..AC) between row 72 and row 73, keep column A (id) as-is
Swap-Rows $ws 72 73 2 29

# Swap data (columns B..AC) between row 101 and row 102, keep column A (id) as-is
Swap-Rows $ws 101 102 2 29
